$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ("aug_0.5"): note the LR used when estimating on the augmented dataset ---
$ws.Range("F3").Value = "Train loss and accuracy are estimated at the end of epoch on augumented dataset. LR 1e-4"

# --- Row 4 ("aug_1.0"): fix the "th" -> "the" typo, fill in the accuracy figures ---
$ws.Range("F4").Value = "1) Added dropout in the fully connected layer`n2) Reduced learning rate from 1e-4 to 5e-5`n3) train loss and accuracy are estimated on the original train dataset, not on the augment dataset"
$ws.Range("D4").Value = 0.83750000000000002
$ws.Range("E4").Value = 0.85833300000000001

# --- Row 5: new "aug_0.5_v2" run, rescaled images to [-1, 1] ---
# Clone row 4's look (fonts / number formats / alignment / wrap) onto row 5 first,
# then overwrite with this run's own data.
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A5").Value = "Baseline3DConvNet"
$ws.Range("B5").Value = "aug_0.5_v2"
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 0.95821699999999999
$ws.Range("E5").Value = 0.92500000000000004
$ws.Range("F5").Value = "1) Max rotation increased from 30 to 135`n2) Rescaled image fro [0, 1] to [-1, 1] range"

# Highlight the valid accuracy of the new run ("Accent4, Lighter 40%" fill).
$ws.Range("E5").Interior.ThemeColor = 8
$ws.Range("E5").Interior.TintAndShade = 0.59999389629810485

$ws.Rows.Item(5).RowHeight = 32

# --- Column F widened to fit the longer notes ---
# (ColumnWidth is in raw character units; the engine adds the usual ~5/6
# character padding when it serialises to the stored XML `width`, so we
# dial the input back by that same amount to land exactly on 48.5.)
$ws.Columns.Item(6).ColumnWidth = 47.666666666666664

# --- Restore the cursor to where the author left it ---
$ws.Range("B6").Select()
